$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, pushing existing rows 81-85 down to 82-86.
$ws.Rows("81:81").Insert()

# Populate the new row 81 with the new weekly record.
$ws.Cells.Item(81, 1).Value = 6
$ws.Cells.Item(81, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 45106
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100112035
$ws.Cells.Item(81, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 580
$ws.Cells.Item(81, 11).Value = 12000
$ws.Cells.Item(81, 12).Value = 13000
$ws.Cells.Item(81, 13).Value = 12448
$ws.Cells.Item(81, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(81, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(81, 16).Value = 830
$ws.Cells.Item(81, 17).Value = 15
$ws.Cells.Item(81, 18).Value = "Hortaliza"
